# Insert a new data row before row 170 (pushes existing rows 170-199 down to 171-200)
# and populate it with the new record's values, matching the source row's format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170; this shifts rows 170..199 down to 171..200
# and should carry the formatting of the surrounding rows (date style on column D, etc.)
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record.
$ws.Cells.Item(170, 1).Value = 10
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170, 3).Value = "La Araucanía"
$ws.Cells.Item(170, 4).Value = 44637
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 100112043
$ws.Cells.Item(170, 7).Value = "Pepino dulce"
$ws.Cells.Item(170, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 250
$ws.Cells.Item(170, 11).Value = 15000
$ws.Cells.Item(170, 12).Value = 17000
$ws.Cells.Item(170, 13).Value = 16200
$ws.Cells.Item(170, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(170, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(170, 16).Value = 900
$ws.Cells.Item(170, 17).Value = 18
$ws.Cells.Item(170, 18).Value = "Hortaliza"
